# Auto-generated Excel COM-interop script
# Applies cached numeric-value updates to the "Ultros_Profits" workbook
# (profit/loss recalculation sheet export) across its 8 item-category sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 28
$ws.Range("H28").Value = 1057.3462
$ws.Range("I28").Value = 886.65
$ws.Range("K28").Value = 886.65
$ws.Range("M28").Value = -401.65
# Row 46
$ws.Range("H46").Value = 783.3333
$ws.Range("I46").Value = 783.3333
$ws.Range("J46").Value = 0
$ws.Range("K46").Value = 2349.9999
$ws.Range("L46").Value = 0
$ws.Range("M46").Value = -2230.9999
$ws.Range("N46").Value = ""
# Row 60
$ws.Range("H60").Value = 783.3333
$ws.Range("I60").Value = 783.3333
$ws.Range("J60").Value = 0
$ws.Range("K60").Value = 2349.9999
$ws.Range("L60").Value = 0
$ws.Range("M60").Value = -1865.9999
$ws.Range("N60").Value = ""
# Row 69
$ws.Range("H69").Value = 37070250
$ws.Range("J69").Value = 37070250
$ws.Range("L69").Value = 111210750
$ws.Range("N69").Value = -111212498
# Row 72
$ws.Range("H72").Value = 37070250
$ws.Range("J72").Value = 37070250
$ws.Range("L72").Value = 333632250
$ws.Range("N72").Value = -333640986
# Row 137
$ws.Range("H137").Value = 5756.2856
$ws.Range("I137").Value = 6850
$ws.Range("K137").Value = 20550
$ws.Range("M137").Value = -18000
$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 33545.273
$ws.Range("I2").Value = 39555.332
$ws.Range("K2").Value = 39555.332
$ws.Range("M2").Value = -39442.332
# Row 30
$ws.Range("H30").Value = 3197.75
$ws.Range("I30").Value = 992
$ws.Range("J30").Value = 3933
$ws.Range("K30").Value = 992
$ws.Range("L30").Value = 3933
$ws.Range("M30").Value = -842
$ws.Range("N30").Value = -4233
# Row 74
$ws.Range("H74").Value = 5565.409
$ws.Range("I74").Value = 4469.3887
$ws.Range("J74").Value = 10497.5
$ws.Range("K74").Value = 4469.3887
$ws.Range("L74").Value = 10497.5
$ws.Range("M74").Value = -3595.3887
$ws.Range("N74").Value = -12245.5
# Row 77
$ws.Range("H77").Value = 5565.409
$ws.Range("I77").Value = 4469.3887
$ws.Range("J77").Value = 10497.5
$ws.Range("K77").Value = 22346.9435
$ws.Range("L77").Value = 52487.5
$ws.Range("M77").Value = -17978.9435
$ws.Range("N77").Value = -61223.5
# Row 88
$ws.Range("H88").Value = 4387258
$ws.Range("I88").Value = 1593.25
$ws.Range("J88").Value = 5556768.5
$ws.Range("K88").Value = 1593.25
$ws.Range("L88").Value = 5556768.5
$ws.Range("M88").Value = -1187.25
$ws.Range("N88").Value = -5557580.5
# Row 91
$ws.Range("H91").Value = 4387258
$ws.Range("I91").Value = 1593.25
$ws.Range("J91").Value = 5556768.5
$ws.Range("K91").Value = 1593.25
$ws.Range("L91").Value = 5556768.5
$ws.Range("M91").Value = -189.25
$ws.Range("N91").Value = -5559576.5
# Row 103
$ws.Range("H103").Value = 27903.5
$ws.Range("J103").Value = 27903.5
$ws.Range("L103").Value = 27903.5
$ws.Range("N103").Value = -30247.5
# Row 116
$ws.Range("H116").Value = 33545.273
$ws.Range("I116").Value = 39555.332
$ws.Range("K116").Value = 39555.332
$ws.Range("M116").Value = -37261.332
$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 33545.273
$ws.Range("I3").Value = 39555.332
$ws.Range("K3").Value = 39555.332
$ws.Range("M3").Value = -39441.332
# Row 76
$ws.Range("H76").Value = 31385.5
$ws.Range("J76").Value = 31385.5
$ws.Range("L76").Value = 31385.5
$ws.Range("N76").Value = -32015.5
# Row 79
$ws.Range("H79").Value = 31385.5
$ws.Range("J79").Value = 31385.5
$ws.Range("L79").Value = 31385.5
$ws.Range("N79").Value = -33569.5
# Row 103
$ws.Range("H103").Value = 21885.666
$ws.Range("J103").Value = 21885.666
$ws.Range("L103").Value = 21885.666
$ws.Range("N103").Value = -24229.666
# Row 107
$ws.Range("H107").Value = 6833.6177
$ws.Range("I107").Value = 6239.5557
$ws.Range("J107").Value = 9125
$ws.Range("K107").Value = 6239.5557
$ws.Range("L107").Value = 9125
$ws.Range("M107").Value = -4319.5557
$ws.Range("N107").Value = -12965
# Row 134
$ws.Range("H134").Value = 2260.7778
$ws.Range("I134").Value = 2286.5386
$ws.Range("K134").Value = 6859.6158
$ws.Range("M134").Value = -4324.6158
$ws = $wb.Worksheets.Item("CRP")
# Row 58
$ws.Range("H58").Value = 2125.6206
$ws.Range("I58").Value = 1300.6471
$ws.Range("K58").Value = 1300.6471
$ws.Range("M58").Value = -1097.6471
# Row 99
$ws.Range("H99").Value = 2622.3333
$ws.Range("I99").Value = 2240.4
$ws.Range("K99").Value = 2240.4
$ws.Range("M99").Value = -742.4000000000001
# Row 126
$ws.Range("H126").Value = 2622.3333
$ws.Range("I126").Value = 2240.4
$ws.Range("K126").Value = 6721.200000000001
$ws.Range("M126").Value = -4251.200000000001
# Row 132
$ws.Range("H132").Value = 3072.8333
$ws.Range("I132").Value = 2545
$ws.Range("J132").Value = 4128.5
$ws.Range("K132").Value = 7635
$ws.Range("L132").Value = 12385.5
$ws.Range("M132").Value = -5105
$ws.Range("N132").Value = -17445.5
# Row 134
$ws.Range("H134").Value = 5591.3887
$ws.Range("I134").Value = 5791.3125
$ws.Range("K134").Value = 17373.9375
$ws.Range("M134").Value = -14838.9375
# Row 136
$ws.Range("H136").Value = 2125.6206
$ws.Range("I136").Value = 1300.6471
$ws.Range("K136").Value = 3901.9413
$ws.Range("M136").Value = -1351.9413
$ws = $wb.Worksheets.Item("CUL")
# Row 14
$ws.Range("H14").Value = 224.2
$ws.Range("I14").Value = 224.2
$ws.Range("K14").Value = 672.5999999999999
$ws.Range("M14").Value = -499.5999999999999
# Row 107
$ws.Range("H107").Value = 1946
$ws.Range("I107").Value = 2416.889
$ws.Range("J107").Value = 1783
$ws.Range("K107").Value = 7250.667
$ws.Range("L107").Value = 5349
$ws.Range("M107").Value = -5330.667
$ws.Range("N107").Value = -9189
# Row 122
$ws.Range("H122").Value = 6061
$ws.Range("I122").Value = 357.6
$ws.Range("J122").Value = 9625.625
$ws.Range("K122").Value = 3218.4
$ws.Range("L122").Value = 86630.625
$ws.Range("M122").Value = -768.4000000000001
$ws.Range("N122").Value = -91530.625
$ws = $wb.Worksheets.Item("GSM")
# Row 43
$ws.Range("H43").Value = 28631.525
$ws.Range("I43").Value = 14833.167
$ws.Range("K43").Value = 14833.167
$ws.Range("M43").Value = -14682.167
# Row 86
$ws.Range("H86").Value = 27571
$ws.Range("J86").Value = 27571
$ws.Range("L86").Value = 27571
$ws.Range("N86").Value = -29943
# Row 89
$ws.Range("H89").Value = 27571
$ws.Range("J89").Value = 27571
$ws.Range("L89").Value = 82713
$ws.Range("N89").Value = -94569
# Row 101
$ws.Range("H101").Value = 32993.2
$ws.Range("J101").Value = 32993.2
$ws.Range("L101").Value = 32993.2
$ws.Range("N101").Value = -39483.2
$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 3590.4
$ws.Range("I7").Value = 2738.25
$ws.Range("K7").Value = 2738.25
$ws.Range("M7").Value = -2626.25
# Row 25
$ws.Range("H25").Value = 46313
$ws.Range("I25").Value = 46313
$ws.Range("K25").Value = 46313
$ws.Range("M25").Value = -46083
# Row 61
$ws.Range("H61").Value = 3475.75
$ws.Range("I61").Value = 3475.75
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 3475.75
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -3273.75
$ws.Range("N61").Value = ""
# Row 99
$ws.Range("H99").Value = 200285
$ws.Range("J99").Value = 200285
$ws.Range("L99").Value = 200285
$ws.Range("N99").Value = -206275
# Row 102
$ws.Range("H102").Value = 200561
$ws.Range("J102").Value = 200561
$ws.Range("L102").Value = 200561
$ws.Range("N102").Value = -207051
# Row 113
$ws.Range("H113").Value = 3475.75
$ws.Range("I113").Value = 3475.75
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 3475.75
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = -1305.75
$ws.Range("N113").Value = ""
# Row 126
$ws.Range("H126").Value = 3590.4
$ws.Range("I126").Value = 2738.25
$ws.Range("K126").Value = 8214.75
$ws.Range("M126").Value = -5744.75
$ws = $wb.Worksheets.Item("WVR")
# Row 8
$ws.Range("H8").Value = 4666
$ws.Range("I8").Value = 4000
$ws.Range("K8").Value = 4000
$ws.Range("M8").Value = -3860
# Row 102
$ws.Range("H102").Value = 62221
$ws.Range("J102").Value = 62221
$ws.Range("L102").Value = 62221
$ws.Range("N102").Value = -68711
# Row 132
$ws.Range("H132").Value = 4524.814
$ws.Range("I132").Value = 3192.4517
$ws.Range("J132").Value = 7966.75
$ws.Range("K132").Value = 9577.355100000001
$ws.Range("L132").Value = 23900.25
$ws.Range("M132").Value = -7047.355100000001
$ws.Range("N132").Value = -28960.25
# Row 136
$ws.Range("H136").Value = 140413.28
$ws.Range("I136").Value = 241662.5
$ws.Range("J136").Value = 5414.3335
$ws.Range("K136").Value = 724987.5
$ws.Range("L136").Value = 16243.0005
$ws.Range("M136").Value = -722437.5
$ws.Range("N136").Value = -21343.0005
